$d = $word.ActiveDocument

# The title paragraph currently reads "Project Genie". We need to remove
# "Genie" while keeping the space that separated it from "Project", and
# that trailing space must end up in its own run (identical formatting
# to the "Project" run), per the target OOXML.

$titlePara = $d.Paragraphs(1)
$titleStart = $titlePara.Range.Start

# Find the word "Genie" within the title paragraph and delete it, leaving
# "Project " (with the trailing space) behind. Use a throwaway duplicate
# range for the Find so the original paragraph position stays valid.
$find = $titlePara.Range.Duplicate.Find
$find.ClearFormatting()
$find.Execute("Genie", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# After deletion, the paragraph text is "Project " (8 characters), starting
# at $titleStart, followed by the paragraph mark. Force the trailing space
# into its own run by toggling a character attribute off and back on --
# this causes the run to split around the space without changing its
# final formatting (it's already bold, so off/on is a no-op on the
# rendered result but forces run segmentation).
$spaceRange = $d.Range($titleStart + 7, $titleStart + 8)
$spaceRange.Font.Bold = $false
$spaceRange.Font.Bold = $true
